$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 369.42856
$ws.Range("I11").Value = 369.42856
$ws.Range("K11").Value = 369.42856
$ws.Range("M11").Value = -229.42856
$ws.Range("H110").Value = 49999.855
$ws.Range("J110").Value = 49999.832
$ws.Range("L110").Value = 49999.832
$ws.Range("N110").Value = -58179.832
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 78
$ws.Range("K2").Value = 78
$ws.Range("M2").Value = 35
$ws.Range("H97").Value = 41668980
$ws.Range("I97").Value = 83336960
$ws.Range("K97").Value = 83336960
$ws.Range("M97").Value = -83336464
$ws.Range("H106").Value = 10000
$ws.Range("J106").Value = 10000
$ws.Range("L106").Value = 10000
$ws.Range("N106").Value = -12524
$ws.Range("H116").Value = 78
$ws.Range("I116").Value = 78
$ws.Range("K116").Value = 78
$ws.Range("M116").Value = 2216
$ws.Range("H119").Value = 75139.39999999999
$ws.Range("J119").Value = 75139.39999999999
$ws.Range("L119").Value = 75139.39999999999
$ws.Range("N119").Value = -84815.39999999999
$ws.Range("H122").Value = 1950
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 78
$ws.Range("I3").Value = 78
$ws.Range("K3").Value = 78
$ws.Range("M3").Value = 36
$ws.Range("H82").Value = 25712.75
$ws.Range("I82").Value = 19024.545
$ws.Range("K82").Value = 19024.545
$ws.Range("M82").Value = -18641.545
$ws.Range("H85").Value = 25712.75
$ws.Range("I85").Value = 19024.545
$ws.Range("K85").Value = 19024.545
$ws.Range("M85").Value = -17698.545
$ws.Range("H86").Value = 6665.5
$ws.Range("I86").Value = 6665.5
$ws.Range("K86").Value = 6665.5
$ws.Range("M86").Value = -5542.5
$ws.Range("H89").Value = 6665.5
$ws.Range("I89").Value = 6665.5
$ws.Range("K89").Value = 33327.5
$ws.Range("M89").Value = -27711.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 23750
$ws.Range("J38").Value = 28000
$ws.Range("L38").Value = 28000
$ws.Range("N38").Value = -28754
$ws.Range("H46").Value = 23750
$ws.Range("J46").Value = 28000
$ws.Range("L46").Value = 28000
$ws.Range("N46").Value = -28422
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 528.125
$ws.Range("I33").Value = 521
$ws.Range("J33").Value = 549.5
$ws.Range("K33").Value = 3126
$ws.Range("L33").Value = 3297
$ws.Range("M33").Value = -2843
$ws.Range("N33").Value = -3863
$ws.Range("H57").Value = 10000
$ws.Range("J57").Value = 10000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31118
$ws.Range("H92").Value = 574.2
$ws.Range("I92").Value = 551.25
$ws.Range("J92").Value = 666
$ws.Range("K92").Value = 1653.75
$ws.Range("L92").Value = 1998
$ws.Range("M92").Value = -405.75
$ws.Range("N92").Value = -4494
$ws.Range("H99").Value = 4749.75
$ws.Range("I99").Value = 4749.75
$ws.Range("K99").Value = 14249.25
$ws.Range("M99").Value = -12003.25
$ws.Range("H103").Value = 372.66666
$ws.Range("J103").Value = 214
$ws.Range("L103").Value = 642
$ws.Range("N103").Value = -2400
$ws.Range("H137").Value = 991.6667
$ws.Range("I137").Value = 991.6667
$ws.Range("K137").Value = 2975.0001
$ws.Range("M137").Value = 2124.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1750
$ws.Range("J6").Value = 1750
$ws.Range("L6").Value = 1750
$ws.Range("N6").Value = -1976
$ws.Range("H16").Value = 1750
$ws.Range("J16").Value = 1750
$ws.Range("L16").Value = 1750
$ws.Range("N16").Value = -2250
$ws.Range("H17").Value = 1163.3334
$ws.Range("J17").Value = 1163.3334
$ws.Range("L17").Value = 1163.3334
$ws.Range("N17").Value = -1499.3334
$ws.Range("H23").Value = 2899.5
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H36").Value = 19999
$ws.Range("J36").Value = 19999
$ws.Range("L36").Value = 19999
$ws.Range("N36").Value = -20969
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 420
$ws.Range("I5").Value = 420
$ws.Range("K5").Value = 420
$ws.Range("M5").Value = -307
$ws.Range("H7").Value = 12830.667
$ws.Range("I7").Value = 10000
$ws.Range("J7").Value = 14246
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 14246
$ws.Range("M7").Value = -9888
$ws.Range("N7").Value = -14470
$ws.Range("H40").Value = 2883.6667
$ws.Range("I40").Value = 2883.6667
$ws.Range("K40").Value = 2883.6667
$ws.Range("M40").Value = -2747.6667
$ws.Range("H122").Value = 3644.6296
$ws.Range("I122").Value = 3108.75
$ws.Range("J122").Value = 4073.3333
$ws.Range("K122").Value = 9326.25
$ws.Range("L122").Value = 12219.9999
$ws.Range("M122").Value = -6876.25
$ws.Range("N122").Value = -17119.9999
$ws.Range("H126").Value = 12830.667
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 14246
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 42738
$ws.Range("M126").Value = -27530
$ws.Range("N126").Value = -47678
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H55").Value = 26499.889
$ws.Range("I55").Value = 4999.5
$ws.Range("J55").Value = 32642.857
$ws.Range("K55").Value = 4999.5
$ws.Range("L55").Value = 32642.857
$ws.Range("M55").Value = -4722.5
$ws.Range("N55").Value = -33196.857
$ws.Range("H117").Value = 33000
$ws.Range("J117").Value = 33000
$ws.Range("L117").Value = 33000
$ws.Range("N117").Value = -42178
$ws.Range("H126").Value = 4492.0713
$ws.Range("I126").Value = 3328.7
$ws.Range("K126").Value = 9986.099999999999
$ws.Range("M126").Value = -7516.099999999999
$ws.Range("H136").Value = 1669.1428
$ws.Range("I136").Value = 1669.1428
$ws.Range("K136").Value = 5007.428400000001
$ws.Range("M136").Value = -2457.428400000001

Write-Host "Applied all edits"